$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the existing row 170 ("1a amarillo" /
# "2a amarillo" pair dated 44427). This pushes the old rows 170-247 down
# to 172-249 and grows the used range to A1:T249.
$ws.Range("A170:A171").EntireRow.Insert()

# Populate the newly inserted row 170 with a fresh "1a amarillo" reading
# (date 44455, $/malla 16 kilos, Región de O'Higgins).
$ws.Range("A170").Value = 11
$ws.Range("B170").Value = "Vega Monumental Concepción"
$ws.Range("C170").Value = "Bíobío"
$ws.Range("D170").Value = 44455
$ws.Range("E170").Value = 8
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100102
$ws.Range("H170").Value = "Cítricos"
$ws.Range("I170").Value = 100102003
$ws.Range("J170").Value = "Limón"
$ws.Range("K170").Value = "Sin especificar"
$ws.Range("L170").Value = "1a amarillo"
$ws.Range("M170").Value = 500
$ws.Range("N170").Value = 6000
$ws.Range("O170").Value = 6000
$ws.Range("P170").Value = 6000
$ws.Range("Q170").Value = "`$/malla 16 kilos"
$ws.Range("R170").Value = "Región de O'Higgins"
$ws.Range("S170").Value = 375
$ws.Range("T170").Value = 16

# Populate the newly inserted row 171 with a fresh "2a amarillo" reading
# (date 44455, $/malla 16 kilos, Región de O'Higgins).
$ws.Range("A171").Value = 11
$ws.Range("B171").Value = "Vega Monumental Concepción"
$ws.Range("C171").Value = "Bíobío"
$ws.Range("D171").Value = 44455
$ws.Range("E171").Value = 8
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100102
$ws.Range("H171").Value = "Cítricos"
$ws.Range("I171").Value = 100102003
$ws.Range("J171").Value = "Limón"
$ws.Range("K171").Value = "Sin especificar"
$ws.Range("L171").Value = "2a amarillo"
$ws.Range("M171").Value = 500
$ws.Range("N171").Value = 5000
$ws.Range("O171").Value = 5000
$ws.Range("P171").Value = 5000
$ws.Range("Q171").Value = "`$/malla 16 kilos"
$ws.Range("R171").Value = "Región de O'Higgins"
$ws.Range("S171").Value = 312
$ws.Range("T171").Value = 16
